# Scheduled runner refresh of market-price snapshots (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H:N)
# across the leve-profit sheets. Values are plain literals scraped from the
# market-board API, so each affected row is updated cell-by-cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3987
$ws.Range("I98").Value = 5230.5
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 5230.5
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -3732.5
$ws.Range("N98").Value = -4496

$ws.Range("H107").Value = 884.0714
$ws.Range("I107").Value = 519.8889
$ws.Range("J107").Value = 1539.6
$ws.Range("K107").Value = 519.8889
$ws.Range("L107").Value = 1539.6
$ws.Range("M107").Value = 1400.1111
$ws.Range("N107").Value = -5379.6

$ws.Range("H112").Value = 1988.7778
$ws.Range("J112").Value = 2271.2856
$ws.Range("L112").Value = 6813.8568
$ws.Range("N112").Value = -9029.856800000001

$ws.Range("H122").Value = 3987
$ws.Range("I122").Value = 5230.5
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 15691.5
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -13241.5
$ws.Range("N122").Value = -9400

$ws.Range("H138").Value = 3518.4092
$ws.Range("I138").Value = 3389.7368
$ws.Range("J138").Value = 4333.3335
$ws.Range("K138").Value = 10169.2104
$ws.Range("L138").Value = 13000.0005
$ws.Range("M138").Value = -5029.2104
$ws.Range("N138").Value = -23280.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2088.0618
$ws.Range("I32").Value = 1767.7534
$ws.Range("J32").Value = 5010.875
$ws.Range("K32").Value = 1767.7534
$ws.Range("L32").Value = 5010.875
$ws.Range("M32").Value = -1480.7534
$ws.Range("N32").Value = -5584.875

$ws.Range("H74").Value = 1690.5
$ws.Range("I74").Value = 1423.4546
$ws.Range("J74").Value = 2110.1428
$ws.Range("K74").Value = 1423.4546
$ws.Range("L74").Value = 2110.1428
$ws.Range("M74").Value = -549.4546
$ws.Range("N74").Value = -3858.1428

$ws.Range("H77").Value = 1690.5
$ws.Range("I77").Value = 1423.4546
$ws.Range("J77").Value = 2110.1428
$ws.Range("K77").Value = 7117.273
$ws.Range("L77").Value = 10550.714
$ws.Range("M77").Value = -2749.273
$ws.Range("N77").Value = -19286.714

$ws.Range("H88").Value = 4603.375
$ws.Range("J88").Value = 4603.375
$ws.Range("L88").Value = 4603.375
$ws.Range("N88").Value = -5415.375

$ws.Range("H91").Value = 4603.375
$ws.Range("J91").Value = 4603.375
$ws.Range("L91").Value = 4603.375
$ws.Range("N91").Value = -7411.375

$ws.Range("H123").Value = 65998.5
$ws.Range("J123").Value = 65998.5
$ws.Range("L123").Value = 65998.5
$ws.Range("N123").Value = -75798.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 226788.78
$ws.Range("I86").Value = 9700
$ws.Range("J86").Value = 335333.16
$ws.Range("K86").Value = 9700
$ws.Range("L86").Value = 335333.16
$ws.Range("M86").Value = -8577
$ws.Range("N86").Value = -337579.16

$ws.Range("H89").Value = 226788.78
$ws.Range("I89").Value = 9700
$ws.Range("J89").Value = 335333.16
$ws.Range("K89").Value = 48500
$ws.Range("L89").Value = 1676665.8
$ws.Range("M89").Value = -42884
$ws.Range("N89").Value = -1687897.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1868.1177
$ws.Range("I31").Value = 1711.9259
$ws.Range("K31").Value = 1711.9259
$ws.Range("M31").Value = -1416.9259

$ws.Range("H34").Value = 1868.1177
$ws.Range("I34").Value = 1711.9259
$ws.Range("K34").Value = 1711.9259
$ws.Range("M34").Value = -1509.9259

$ws.Range("H58").Value = 1036352.5
$ws.Range("I58").Value = 1403474.1
$ws.Range("J58").Value = 1737.2727
$ws.Range("K58").Value = 1403474.1
$ws.Range("L58").Value = 1737.2727
$ws.Range("M58").Value = -1403271.1
$ws.Range("N58").Value = -2143.2727

$ws.Range("H132").Value = 1285.0278
$ws.Range("I132").Value = 863.13336
$ws.Range("K132").Value = 2589.40008
$ws.Range("M132").Value = -59.40008000000034

$ws.Range("H136").Value = 1036352.5
$ws.Range("I136").Value = 1403474.1
$ws.Range("J136").Value = 1737.2727
$ws.Range("K136").Value = 4210422.300000001
$ws.Range("L136").Value = 5211.8181
$ws.Range("M136").Value = -4207872.300000001
$ws.Range("N136").Value = -10311.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 722.8889
$ws.Range("J5").Value = 868.05884
$ws.Range("L5").Value = 2604.17652
$ws.Range("N5").Value = -2828.17652

$ws.Range("H26").Value = 607.4
$ws.Range("I26").Value = 782
$ws.Range("K26").Value = 2346
$ws.Range("M26").Value = -2058

$ws.Range("H68").Value = 777.63635
$ws.Range("I68").Value = 782.6667
$ws.Range("K68").Value = 2348.0001
$ws.Range("M68").Value = -1537.0001

$ws.Range("H71").Value = 777.63635
$ws.Range("I71").Value = 782.6667
$ws.Range("K71").Value = 7044.0003
$ws.Range("M71").Value = -2988.0003

$ws.Range("H103").Value = 2106.652
$ws.Range("I103").Value = 2711
$ws.Range("J103").Value = 1938.7778
$ws.Range("K103").Value = 8133
$ws.Range("L103").Value = 5816.3334
$ws.Range("M103").Value = -7254
$ws.Range("N103").Value = -7574.3334

$ws.Range("H109").Value = 2443.5
$ws.Range("I109").Value = 1222.375
$ws.Range("J109").Value = 4071.6667
$ws.Range("K109").Value = 3667.125
$ws.Range("L109").Value = 12215.0001
$ws.Range("M109").Value = -2627.125
$ws.Range("N109").Value = -14295.0001

$ws.Range("H114").Value = 3390.8
$ws.Range("I114").Value = 652
$ws.Range("J114").Value = 7499
$ws.Range("K114").Value = 1956
$ws.Range("L114").Value = 22497
$ws.Range("M114").Value = 1298
$ws.Range("N114").Value = -29005

$ws.Range("H131").Value = 13223.8
$ws.Range("I131").Value = 842
$ws.Range("J131").Value = 14255.616
$ws.Range("K131").Value = 2526
$ws.Range("L131").Value = 42766.848
$ws.Range("M131").Value = 2514
$ws.Range("N131").Value = -52846.848

$ws.Range("H135").Value = 722.8889
$ws.Range("J135").Value = 868.05884
$ws.Range("L135").Value = 7812.52956
$ws.Range("N135").Value = -12882.52956

$ws.Range("H139").Value = 14889.125
$ws.Range("I139").Value = 18852.166
$ws.Range("K139").Value = 56556.49800000001
$ws.Range("M139").Value = -51416.49800000001

$ws.Range("H140").Value = 2260.2632
$ws.Range("I140").Value = 1423.625
$ws.Range("J140").Value = 2868.7273
$ws.Range("K140").Value = 4270.875
$ws.Range("L140").Value = 8606.1819
$ws.Range("M140").Value = 909.125
$ws.Range("N140").Value = -18966.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5146660
$ws.Range("I11").Value = 6174462
$ws.Range("K11").Value = 6174462
$ws.Range("M11").Value = -6174323

$ws.Range("H70").Value = 4300
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4300
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4300
$ws.Range("N70").Value = -4840
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 4300
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4300
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4300
$ws.Range("N73").Value = -6172
$ws.Range("M73").ClearContents()

$ws.Range("H122").Value = 1364.2916
$ws.Range("I122").Value = 1107.7333
$ws.Range("J122").Value = 1791.8889
$ws.Range("K122").Value = 3323.199900000001
$ws.Range("L122").Value = 5375.6667
$ws.Range("M122").Value = -873.1999000000005
$ws.Range("N122").Value = -10275.6667

$ws.Range("H132").Value = 786852.9
$ws.Range("I132").Value = 1166840.5
$ws.Range("K132").Value = 3500521.5
$ws.Range("M132").Value = -3497991.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2875.75
$ws.Range("I7").Value = 1834.3334
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 1834.3334
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -1722.3334
$ws.Range("N7").Value = -6224

$ws.Range("H61").Value = 2235.353
$ws.Range("I61").Value = 2227.6365
$ws.Range("J61").Value = 2249.5
$ws.Range("K61").Value = 2227.6365
$ws.Range("L61").Value = 2249.5
$ws.Range("M61").Value = -2025.6365
$ws.Range("N61").Value = -2653.5

$ws.Range("H68").Value = 2466.6667
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 2466.6667
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488

$ws.Range("H82").Value = 1303.1538
$ws.Range("I82").Value = 1110.2222
$ws.Range("K82").Value = 1110.2222
$ws.Range("M82").Value = -749.2221999999999

$ws.Range("H85").Value = 1303.1538
$ws.Range("I85").Value = 1110.2222
$ws.Range("K85").Value = 1110.2222
$ws.Range("M85").Value = 137.7778000000001

$ws.Range("H113").Value = 2235.353
$ws.Range("I113").Value = 2227.6365
$ws.Range("J113").Value = 2249.5
$ws.Range("K113").Value = 2227.6365
$ws.Range("L113").Value = 2249.5
$ws.Range("M113").Value = -57.63650000000007
$ws.Range("N113").Value = -6589.5

$ws.Range("H126").Value = 2875.75
$ws.Range("I126").Value = 1834.3334
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 5503.0002
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -3033.0002
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 26457540
$ws.Range("I136").Value = 42737480
$ws.Range("K136").Value = 128212440
$ws.Range("M136").Value = -128209890
